$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.574.99'
$ws.Range("E2").Value = '  +3.99%  '
# Row 3
$ws.Range("D3").Value = '2.259.68'
$ws.Range("E3").Value = '  +1.07%  '
# Row 4
$ws.Range("E4").Value = '  +0.07%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.55'
$ws.Range("E5").Value = '  -0.51%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  +0.60%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.24'
$ws.Range("E7").Value = '  -0.95%  '
# Row 8
$ws.Range("E8").Value = '  +0.09%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.420'
$ws.Range("E9").Value = '  +4.18%  '
# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.89'
$ws.Range("E10").Value = '  -2.41%  '
# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("E11").Value = '  +4.73%  '
# Row 12
$ws.Range("E12").Value = '  +0.47%  '
# Row 13
$ws.Range("D13").Value = '2.597.28'
$ws.Range("E13").Value = '  +1.38%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.57'
$ws.Range("E14").Value = '  -0.84%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.58'
$ws.Range("E15").Value = '  +6.55%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.76'
$ws.Range("E16").Value = '  +2.74%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.808'
$ws.Range("E17").Value = '  +0.63%  '
# Row 18
$ws.Range("D18").Value = '2.265.23'
$ws.Range("E18").Value = '  +0.86%  '
# Row 19
$ws.Range("D19").Value = '43.446.60'
$ws.Range("E19").Value = '  +3.87%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0931'
$ws.Range("E20").Value = '  +3.94%  '
# Row 21
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.77'
$ws.Range("E21").Value = '  +0.92%  '
# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  +2.44%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.18'
$ws.Range("E23").Value = '  +0.66%  '
# Row 24
$ws.Range("E24").Value = '  -0.10%  '
# Row 25
$ws.Range("E25").Value = '  +6.01%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'
$ws.Range("E26").Value = '  +1.92%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.80'
$ws.Range("E27").Value = '  +1.39%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '170.54'
$ws.Range("E28").Value = '  +2.22%  '
# Row 29
$ws.Range("E29").Value = '  -1.85%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.45'
$ws.Range("E30").Value = '  +2.35%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("E31").Value = '  +1.88%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.65'
$ws.Range("E32").Value = '  +0.83%  '
# Row 33
$ws.Range("E33").Value = '  -0.25%  '
# Row 34
$ws.Range("E34").Value = '  +0.32%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.76'
$ws.Range("E35").Value = '  +1.38%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0657'
$ws.Range("E36").Value = '  +3.21%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.43'
$ws.Range("E37").Value = '  -3.62%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.39'
$ws.Range("E38").Value = '  +0.38%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.58'
$ws.Range("E39").Value = '  -2.62%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0249'
$ws.Range("E40").Value = '  +3.67%  '
# Row 41
$ws.Range("E41").Value = '  +0.23%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000230'
$ws.Range("E42").Value = '  -10.75%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.69'
$ws.Range("E43").Value = '  +1.23%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0988'
$ws.Range("E44").Value = '  +0.69%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.52'
$ws.Range("E45").Value = '  -6.42%  '
# Row 46
$ws.Range("E46").Value = '  -1.10%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.76'
$ws.Range("E47").Value = '  -1.30%  '
# Row 48
$ws.Range("D48").Value = '1.468.65'
$ws.Range("E48").Value = '  -0.65%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.56'
$ws.Range("E49").Value = '  -0.07%  '
# Row 50
$ws.Range("E50").Value = '  +0.13%  '
# Row 51
$ws.Range("E51").Value = '  +7.13%  '
